$wb = $excel.ActiveWorkbook

# Update "想去人数" (interested-count) figures for the four events that
# appear both on the "展览" sheet and on the consolidated "全部类型" sheet.
$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F3").Value = 106
    $ws.Range("F4").Value = 129
    $ws.Range("F5").Value = 2848
    $ws.Range("F6").Value = 282
}
